{"js": "// Replace the date line and every \"AA\u00d7BB=\" problem in the table with the\n// new values from the commit. Each old string is unique in the document,\n// so a plain (case-sensitive, non-wildcard) search/replace per pair is\n// sufficient and safe - no two replacements collide with each other's\n// before/after text.\nconst replacements = [\n    [\"2024-08-29 Thursday\", \"2024-08-30 Friday\"],\n    [\"87\u00d732=\", \"23\u00d789=\"],\n    [\"26\u00d736=\", \"87\u00d720=\"],\n    [\"16\u00d779=\", \"87\u00d759=\"],\n    [\"34\u00d794=\", \"55\u00d729=\"],\n    [\"54\u00d760=\", \"65\u00d747=\"],\n    [\"95\u00d779=\", \"16\u00d757=\"],\n    [\"80\u00d773=\", \"15\u00d777=\"],\n    [\"12\u00d764=\", \"55\u00d769=\"],\n    [\"71\u00d759=\", \"27\u00d718=\"],\n    [\"97\u00d742=\", \"31\u00d780=\"],\n    [\"37\u00d799=\", \"78\u00d797=\"],\n    [\"33\u00d769=\", \"27\u00d719=\"],\n    [\"26\u00d724=\", \"47\u00d735=\"],\n    [\"59\u00d747=\", \"54\u00d770=\"],\n    [\"39\u00d774=\", \"91\u00d760=\"],\n    [\"76\u00d761=\", \"90\u00d729=\"],\n    [\"80\u00d779=\", \"63\u00d796=\"],\n    [\"54\u00d769=\", \"26\u00d727=\"],\n    [\"91\u00d748=\", \"66\u00d750=\"],\n    [\"19\u00d725=\", \"62\u00d781=\"],\n    [\"81\u00d734=\", \"93\u00d736=\"],\n    [\"41\u00d714=\", \"34\u00d755=\"],\n    [\"77\u00d756=\", \"14\u00d713=\"],\n    [\"54\u00d750=\", \"74\u00d767=\"],\n    [\"58\u00d750=\", \"42\u00d750=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (let i = 0; i < results.items.length; i++) {\n        results.items[i].insertText(newText, Word.InsertLocation.replace);\n    }\n    await context.sync();\n}\n", "ps1": "# Replace the date line and every \"AA\u00d7BB=\" problem in the table with the\n# new values from the commit. Each old string is unique in the document,\n# so a plain (case-sensitive, non-wildcard) Find/Replace per pair is\n# sufficient and safe - no two replacements collide with each other's\n# before/after text.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-08-29 Thursday\", \"2024-08-30 Friday\"),\n    @(\"87\u00d732=\", \"23\u00d789=\"),\n    @(\"26\u00d736=\", \"87\u00d720=\"),\n    @(\"16\u00d779=\", \"87\u00d759=\"),\n    @(\"34\u00d794=\", \"55\u00d729=\"),\n    @(\"54\u00d760=\", \"65\u00d747=\"),\n    @(\"95\u00d779=\", \"16\u00d757=\"),\n    @(\"80\u00d773=\", \"15\u00d777=\"),\n    @(\"12\u00d764=\", \"55\u00d769=\"),\n    @(\"71\u00d759=\", \"27\u00d718=\"),\n    @(\"97\u00d742=\", \"31\u00d780=\"),\n    @(\"37\u00d799=\", \"78\u00d797=\"),\n    @(\"33\u00d769=\", \"27\u00d719=\"),\n    @(\"26\u00d724=\", \"47\u00d735=\"),\n    @(\"59\u00d747=\", \"54\u00d770=\"),\n    @(\"39\u00d774=\", \"91\u00d760=\"),\n    @(\"76\u00d761=\", \"90\u00d729=\"),\n    @(\"80\u00d779=\", \"63\u00d796=\"),\n    @(\"54\u00d769=\", \"26\u00d727=\"),\n    @(\"91\u00d748=\", \"66\u00d750=\"),\n    @(\"19\u00d725=\", \"62\u00d781=\"),\n    @(\"81\u00d734=\", \"93\u00d736=\"),\n    @(\"41\u00d714=\", \"34\u00d755=\"),\n    @(\"77\u00d756=\", \"14\u00d713=\"),\n    @(\"54\u00d750=\", \"74\u00d767=\"),\n    @(\"58\u00d750=\", \"42\u00d750=\")\n)\n\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n#              MatchSoundsLike, MatchAllWordForms, Forward, Wrap,\n#              Format, ReplaceWith, Replace)\n# wdFindContinue = 1, wdReplaceAll = 2\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
